$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, shifting existing rows 39-54 down to 40-55.
$ws.Rows(39).Insert()

# Populate the newly inserted row 39 with the new observation.
$ws.Range("A39").Value = 1
$ws.Range("B39").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C39").Value = "Arica y Parinacota"
$ws.Range("D39").Value = 44876
$ws.Range("E39").Value = 15
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = "Frutos de hueso (carozo)"
$ws.Range("I39").Value = 100103004
$ws.Range("J39").Value = "Durazno"
$ws.Range("K39").Value = "Florida King"
$ws.Range("L39").Value = "Segunda"
$ws.Range("M39").Value = 300
$ws.Range("N39").Value = 22000
$ws.Range("O39").Value = 23000
$ws.Range("P39").Value = 22500
$ws.Range("Q39").Value = "$/bandeja 10 kilos granel"
$ws.Range("R39").Value = "Región de Coquimbo"
$ws.Range("S39").Value = 2250
$ws.Range("T39").Value = 10
